# Applies the "suivi_projet" update: new activities on 2025-06-03..06-05 (3 Jun
# to 5 Jun), a new blank-date sub-row inserted for 2025-06-05, and small
# cosmetic sheet-view adjustments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 2 (03/06) grows to a two-line entry -> explicit row height.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Row 28 (the continuation row for 03/06, merged under A27) gets a new
#    duration value (5h instead of 3h).
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = "5h"

# ---------------------------------------------------------------------------
# 3. Insert one new row above the current row 31 (the 06/06 row). This makes
#    room for a second sub-task on 05/06 while pushing 06/06..20/06 down by
#    one row (old row31 -> row32, ..., old row45 -> row46).
# ---------------------------------------------------------------------------
$ws.Rows(31).Insert()

# ---------------------------------------------------------------------------
# 4. Fill in the now-existing row 29 (04/06) with its task.
# ---------------------------------------------------------------------------
$ws.Range("B29").Value = "3h"
$ws.Range("C29").Value = "Développement application QtCreator avec l'interface + formation trame avec le checksum + creation BDD"
$ws.Rows(29).RowHeight = 60

# ---------------------------------------------------------------------------
# 5. Fill in row 30 (05/06, first sub-task).
# ---------------------------------------------------------------------------
$ws.Range("B30").Value = "1h"
$ws.Range("C30").Value = "Développement application avec intégration de la BDD à l'afficheur"
$ws.Rows(30).RowHeight = 30

# ---------------------------------------------------------------------------
# 6. Fill in the newly inserted row 31 (05/06, second sub-task, blank date
#    cell merged with A30).
# ---------------------------------------------------------------------------
$ws.Range("B31").Value = "2h"
$ws.Range("C31").Value = "Gérer la création et la suppression d'indice de la BDD depuis l'application"
$ws.Rows(31).RowHeight = 30

# ---------------------------------------------------------------------------
# 7. Merge the date cells that now span two rows each.
# ---------------------------------------------------------------------------
$ws.Range("A27:A28").Merge()
$ws.Range("A30:A31").Merge()

# ---------------------------------------------------------------------------
# 8. Restore/align formatting for the date column cells touched above so
#    they keep matching the rest of the "date" column.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 9. Conditional-formatting "today" rule no longer needs to cover the blank
#    merged date cells (A28, A31) individually - keep it excluding them, same
#    as it already excluded A28 before this edit.
# ---------------------------------------------------------------------------
$ws.Range("A1:A9,A12,A17:A18,A20:A22,A27,A29:A30,A32:A1048576").FormatConditions(1).Delete()
$fc = $ws.Range("A1:A9,A12,A17:A18,A20:A22,A27,A29:A30,A32:A1048576").FormatConditions.Add(2, 0, "FLOOR(A1,1)=TODAY()")

# ---------------------------------------------------------------------------
# 10. Sheet view: keep the selection where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("C31").Select()
